$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.5416
$ws.Range("B7").Value = 4.900800000000003
$ws.Range("A10").Value = -21.8585
$ws.Range("A12").Value = -21.5715
$ws.Range("B15").Value = 4.562899999999996
$ws.Range("A18").Value = -22.14610000000001
$ws.Range("E18").Value = 18.11390000000002
$ws.Range("E19").Value = 16.2458
$ws.Range("B20").Value = 9.323999999999993
$ws.Range("E27").Value = 16.46789999999999
$ws.Range("B29").Value = 4.911900000000005
$ws.Range("B30").Value = 4.848900000000002
$ws.Range("B31").Value = 5.792199999999999
$ws.Range("A37").Value = -19.31959999999999
$ws.Range("B40").Value = 9.501999999999992
$ws.Range("E42").Value = 16.37329999999999
$ws.Range("E44").Value = 16.35389999999999
$ws.Range("E47").Value = 16.43059999999999
$ws.Range("A55").Value = -21.7849
$ws.Range("E58").Value = 16.32030000000002
$ws.Range("A68").Value = -21.50770000000001
$ws.Range("B68").Value = 4.551900000000001
$ws.Range("E73").Value = 17.40870000000001
$ws.Range("B76").Value = 5.631499999999997
$ws.Range("A77").Value = -20.51739999999999
$ws.Range("A78").Value = -20.03439999999998
$ws.Range("B87").Value = 4.589499999999993
$ws.Range("B88").Value = 4.644499999999997
$ws.Range("E95").Value = 18.25160000000002
$ws.Range("B96").Value = 5.459400000000006
$ws.Range("B98").Value = 5.946299999999999
$ws.Range("B101").Value = 9.625999999999998
$ws.Range("E101").Value = 16.60610000000001
$ws.Range("B102").Value = 8.375500000000004
